$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row based on the sheet's used range.
$lastRow = $ws.UsedRange.Rows.Count

# Swap the contents of columns C (codeforiati:group-name) and D (codeforiati:group-code)
# for every row, including the header row, so that column C becomes the code
# and column D becomes the name.
for ($r = 1; $r -le $lastRow; $r++) {
    $cVal = $ws.Cells.Item($r, 3).Value2
    $dVal = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 3).Value = $dVal
    $ws.Cells.Item($r, 4).Value = $cVal
}
